# Auto-generated edit script applying the diff from the commit
# "chore: update Sheets via scheduled runner"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 1499.875
$ws.Range("I51").Value = 1499
$ws.Range("K51").Value = 1499
$ws.Range("M51").Value = -1015
$ws.Range("H61").Value = 996.25
$ws.Range("I61").Value = 996.25
$ws.Range("K61").Value = 2988.75
$ws.Range("M61").Value = -2816.75
$ws.Range("H74").Value = 7928.5713
$ws.Range("I74").Value = 7928.5713
$ws.Range("K74").Value = 7928.5713
$ws.Range("M74").Value = -6992.5713
$ws.Range("H77").Value = 7928.5713
$ws.Range("I77").Value = 7928.5713
$ws.Range("K77").Value = 39642.85649999999
$ws.Range("M77").Value = -34962.85649999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 5550.0835
$ws.Range("I102").Value = 5600.091
$ws.Range("K102").Value = 5600.091
$ws.Range("M102").Value = -3978.091
$ws.Range("H122").Value = 2463.1428
$ws.Range("I122").Value = 2658.7
$ws.Range("K122").Value = 7976.099999999999
$ws.Range("M122").Value = -5526.099999999999
$ws.Range("H132").Value = 1928.7
$ws.Range("I132").Value = 1138.4667
$ws.Range("K132").Value = 3415.4001
$ws.Range("M132").Value = -885.4000999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8997.799999999999
$ws.Range("I86").Value = 3995
$ws.Range("J86").Value = 12333
$ws.Range("K86").Value = 3995
$ws.Range("L86").Value = 12333
$ws.Range("M86").Value = -2872
$ws.Range("N86").Value = -14579
$ws.Range("H89").Value = 8997.799999999999
$ws.Range("I89").Value = 3995
$ws.Range("J89").Value = 12333
$ws.Range("K89").Value = 19975
$ws.Range("L89").Value = 61665
$ws.Range("M89").Value = -14359
$ws.Range("N89").Value = -72897
$ws.Range("H94").Value = 3256.4443
$ws.Range("I94").Value = 2385
$ws.Range("J94").Value = 4999.3335
$ws.Range("K94").Value = 2385
$ws.Range("L94").Value = 4999.3335
$ws.Range("M94").Value = -1934
$ws.Range("N94").Value = -5901.3335
$ws.Range("H99").Value = 2193.75
$ws.Range("I99").Value = 1658.3334
$ws.Range("J99").Value = 3800
$ws.Range("K99").Value = 1658.3334
$ws.Range("L99").Value = 3800
$ws.Range("M99").Value = -160.3334
$ws.Range("N99").Value = -6796
$ws.Range("H105").Value = 4369.8
$ws.Range("I105").Value = 4369.8
$ws.Range("K105").Value = 4369.8
$ws.Range("M105").Value = -2622.8
$ws.Range("H134").Value = 2096.6
$ws.Range("I134").Value = 1981
$ws.Range("J134").Value = 2366.3333
$ws.Range("K134").Value = 5943
$ws.Range("L134").Value = 7098.999899999999
$ws.Range("M134").Value = -3408
$ws.Range("N134").Value = -12168.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 671.8
$ws.Range("I16").Value = 586.5
$ws.Range("K16").Value = 586.5
$ws.Range("M16").Value = -299.5
$ws.Range("H31").Value = 1793.5625
$ws.Range("I31").Value = 1378
$ws.Range("J31").Value = 2209.125
$ws.Range("K31").Value = 1378
$ws.Range("L31").Value = 2209.125
$ws.Range("M31").Value = -1083
$ws.Range("N31").Value = -2799.125
$ws.Range("H34").Value = 1793.5625
$ws.Range("I34").Value = 1378
$ws.Range("J34").Value = 2209.125
$ws.Range("K34").Value = 1378
$ws.Range("L34").Value = 2209.125
$ws.Range("M34").Value = -1176
$ws.Range("N34").Value = -2613.125
$ws.Range("H113").Value = 671.8
$ws.Range("I113").Value = 586.5
$ws.Range("K113").Value = 586.5
$ws.Range("M113").Value = 1583.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1501740.4
$ws.Range("I8").Value = 1501740.4
$ws.Range("K8").Value = 4505221.199999999
$ws.Range("M8").Value = -4505082.199999999
$ws.Range("H38").Value = 241.25
$ws.Range("J38").Value = 430.25
$ws.Range("L38").Value = 1290.75
$ws.Range("N38").Value = -1984.75
$ws.Range("H132").Value = 1263.6666
$ws.Range("I132").Value = 1263.6666
$ws.Range("K132").Value = 11372.9994
$ws.Range("M132").Value = -8842.999400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 9592.25
$ws.Range("I80").Value = 2789.6667
$ws.Range("J80").Value = 30000
$ws.Range("K80").Value = 2789.6667
$ws.Range("L80").Value = 30000
$ws.Range("M80").Value = -1791.6667
$ws.Range("N80").Value = -31996
$ws.Range("H83").Value = 9592.25
$ws.Range("I83").Value = 2789.6667
$ws.Range("J83").Value = 30000
$ws.Range("K83").Value = 13948.3335
$ws.Range("L83").Value = 150000
$ws.Range("M83").Value = -8956.333500000001
$ws.Range("N83").Value = -159984
$ws.Range("H97").Value = 7671.5
$ws.Range("I97").Value = 9992.666999999999
$ws.Range("J97").Value = 708
$ws.Range("K97").Value = 9992.666999999999
$ws.Range("L97").Value = 708
$ws.Range("M97").Value = -9496.666999999999
$ws.Range("N97").Value = -1700
$ws.Range("H122").Value = 3666.75
$ws.Range("I122").Value = 3741
$ws.Range("J122").Value = 3444
$ws.Range("K122").Value = 11223
$ws.Range("L122").Value = 10332
$ws.Range("M122").Value = -8773
$ws.Range("N122").Value = -15232

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6811
$ws.Range("I22").Value = 4685
$ws.Range("K22").Value = 4685
$ws.Range("M22").Value = -4390
$ws.Range("H27").Value = 6811
$ws.Range("I27").Value = 4685
$ws.Range("K27").Value = 4685
$ws.Range("M27").Value = -4578
$ws.Range("H40").Value = 3459.6
$ws.Range("I40").Value = 3459.6
$ws.Range("K40").Value = 3459.6
$ws.Range("M40").Value = -3323.6
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H93").Value = 4501.5
$ws.Range("I93").Value = 4501.5
$ws.Range("K93").Value = 4501.5
$ws.Range("M93").Value = -3253.5
$ws.Range("H100").Value = 4113.1665
$ws.Range("I100").Value = 4113.1665
$ws.Range("K100").Value = 4113.1665
$ws.Range("M100").Value = -3572.1665
$ws.Range("H132").Value = 5316.1665
$ws.Range("I132").Value = 4000
$ws.Range("K132").Value = 12000
$ws.Range("M132").Value = -9470
$ws.Range("H136").Value = 2687
$ws.Range("I136").Value = 2687
$ws.Range("K136").Value = 8061
$ws.Range("M136").Value = -5511

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2025.0667
$ws.Range("I100").Value = 1597.125
$ws.Range("J100").Value = 2514.1428
$ws.Range("K100").Value = 3194.25
$ws.Range("L100").Value = 5028.2856
$ws.Range("M100").Value = -2653.25
$ws.Range("N100").Value = -6110.2856
$ws.Range("H132").Value = 2267
$ws.Range("I132").Value = 1626.2
$ws.Range("K132").Value = 4878.6
$ws.Range("M132").Value = -2348.6
